$wb = $excel.ActiveWorkbook
$typ = $wb.Worksheets.Item("Typography")
$trans = $wb.Worksheets.Item("Translation")

# --- Typography sheet: new rows 9-12 (displayMeas, displayLabel, displaylabelBold, displayMeasValue) ---
# Row 9
$typ.Cells.Item(9, 2).Value = "displayMeas"
$typ.Cells.Item(9, 3).Value = "verdana.ttf"
$typ.Cells.Item(9, 4).Value = 22
$typ.Cells.Item(9, 5).Value = 4
$typ.Cells.Item(9, 6).Value = "?"
$typ.Cells.Item(9, 9).Value = "a-z,A-Z,0-9"

# Row 10
$typ.Cells.Item(10, 2).Value = "displayLabel"
$typ.Cells.Item(10, 3).Value = "verdana.ttf"
$typ.Cells.Item(10, 4).Value = 15
$typ.Cells.Item(10, 5).Value = 4
$typ.Cells.Item(10, 6).Value = "?"

# Row 11
$typ.Cells.Item(11, 2).Value = "displaylabelBold"
$typ.Cells.Item(11, 3).Value = "verdanab.ttf"
$typ.Cells.Item(11, 4).Value = 15
$typ.Cells.Item(11, 5).Value = 4
$typ.Cells.Item(11, 6).Value = "?"

# Row 12
$typ.Cells.Item(12, 2).Value = "displayMeasValue"
$typ.Cells.Item(12, 3).Value = "verdanab.ttf"
$typ.Cells.Item(12, 4).Value = 18
$typ.Cells.Item(12, 5).Value = 4
$typ.Cells.Item(12, 6).Value = "?"

# --- Translation sheet: new rows 75-127 (time-mode display measure panel) ---
# Row 75
$trans.Cells.Item(75, 2).Value = "SingleUseId90"
$trans.Cells.Item(75, 3).Value = "displayMeas"
$trans.Cells.Item(75, 4).Value = "Left"
$trans.Cells.Item(75, 5).Value = "LTR"
$trans.Cells.Item(75, 6).Value = "TI 1 (IN<value>"

# Row 76
$trans.Cells.Item(76, 2).Value = "SingleUseId91"
$trans.Cells.Item(76, 3).Value = "displayMeas"
$trans.Cells.Item(76, 4).Value = "Left"
$trans.Cells.Item(76, 5).Value = "LTR"
$trans.Cells.Item(76, 6).Value = "IN<value>)"

# Row 77
$trans.Cells.Item(77, 2).Value = "SingleUseId92"
$trans.Cells.Item(77, 3).Value = "displayMeas"
$trans.Cells.Item(77, 4).Value = "Left"
$trans.Cells.Item(77, 5).Value = "LTR"
$trans.Cells.Item(77, 6).Value = "TI 2 (IN<value>"

# Row 78
$trans.Cells.Item(78, 2).Value = "SingleUseId93"
$trans.Cells.Item(78, 3).Value = "displayMeas"
$trans.Cells.Item(78, 4).Value = "Left"
$trans.Cells.Item(78, 5).Value = "LTR"
$trans.Cells.Item(78, 6).Value = "IN<value>)"

# Row 79
$trans.Cells.Item(79, 2).Value = "SingleUseId94"
$trans.Cells.Item(79, 3).Value = "displayMeas"
$trans.Cells.Item(79, 4).Value = "Left"
$trans.Cells.Item(79, 5).Value = "LTR"
$trans.Cells.Item(79, 6).Value = "TI 3 (IN<value>"

# Row 80
$trans.Cells.Item(80, 2).Value = "SingleUseId95"
$trans.Cells.Item(80, 3).Value = "displayMeas"
$trans.Cells.Item(80, 4).Value = "Left"
$trans.Cells.Item(80, 5).Value = "LTR"
$trans.Cells.Item(80, 6).Value = "IN<value>)"

# Row 81
$trans.Cells.Item(81, 2).Value = "SingleUseId96"
$trans.Cells.Item(81, 3).Value = "displayMeas"
$trans.Cells.Item(81, 4).Value = "Left"
$trans.Cells.Item(81, 5).Value = "LTR"
$trans.Cells.Item(81, 6).Value = "TI 4 (IN<value>"

# Row 82
$trans.Cells.Item(82, 2).Value = "SingleUseId97"
$trans.Cells.Item(82, 3).Value = "displayMeas"
$trans.Cells.Item(82, 4).Value = "Left"
$trans.Cells.Item(82, 5).Value = "LTR"
$trans.Cells.Item(82, 6).Value = "IN<value>)"

# Row 83
$trans.Cells.Item(83, 2).Value = "SingleUseId98"
$trans.Cells.Item(83, 3).Value = "displayMeas"
$trans.Cells.Item(83, 4).Value = "Left"
$trans.Cells.Item(83, 5).Value = "LTR"
$trans.Cells.Item(83, 6).Value = "TI 5 (IN<value>"

# Row 84
$trans.Cells.Item(84, 2).Value = "SingleUseId99"
$trans.Cells.Item(84, 3).Value = "displayMeas"
$trans.Cells.Item(84, 4).Value = "Left"
$trans.Cells.Item(84, 5).Value = "LTR"
$trans.Cells.Item(84, 6).Value = "IN<value>)"

# Row 85
$trans.Cells.Item(85, 2).Value = "SingleUseId100"
$trans.Cells.Item(85, 3).Value = "displayMeas"
$trans.Cells.Item(85, 4).Value = "Left"
$trans.Cells.Item(85, 5).Value = "LTR"
$trans.Cells.Item(85, 6).Value = "TI 6 (IN<value>"

# Row 86
$trans.Cells.Item(86, 2).Value = "SingleUseId101"
$trans.Cells.Item(86, 3).Value = "displayMeas"
$trans.Cells.Item(86, 4).Value = "Left"
$trans.Cells.Item(86, 5).Value = "LTR"
$trans.Cells.Item(86, 6).Value = "IN<value>)"

# Row 87
$trans.Cells.Item(87, 2).Value = "SingleUseId102"
$trans.Cells.Item(87, 3).Value = "displayMeas"
$trans.Cells.Item(87, 4).Value = "Left"
$trans.Cells.Item(87, 5).Value = "LTR"
$trans.Cells.Item(87, 6).Value = "TI 7 (IN<value>"

# Row 88
$trans.Cells.Item(88, 2).Value = "SingleUseId103"
$trans.Cells.Item(88, 3).Value = "displayMeas"
$trans.Cells.Item(88, 4).Value = "Left"
$trans.Cells.Item(88, 5).Value = "LTR"
$trans.Cells.Item(88, 6).Value = "IN<value>)"

# Row 89
$trans.Cells.Item(89, 2).Value = "SingleUseId104"
$trans.Cells.Item(89, 3).Value = "displayLabel"
$trans.Cells.Item(89, 4).Value = "Left"
$trans.Cells.Item(89, 5).Value = "LTR"
$trans.Cells.Item(89, 6).Value = "IN<value>: <value>"

# Row 90
$trans.Cells.Item(90, 2).Value = "SingleUseId105"
$trans.Cells.Item(90, 3).Value = "displaylabelBold"
$trans.Cells.Item(90, 4).Value = "Left"
$trans.Cells.Item(90, 5).Value = "LTR"
$trans.Cells.Item(90, 6).Value = "Stamps "

# Row 91
$trans.Cells.Item(91, 2).Value = "SingleUseId106"
$trans.Cells.Item(91, 3).Value = "displaylabelBold"
$trans.Cells.Item(91, 4).Value = "Left"
$trans.Cells.Item(91, 5).Value = "LTR"
$trans.Cells.Item(91, 6).Value = "Sample"

# Row 92
$trans.Cells.Item(92, 2).Value = "SingleUseId107"
$trans.Cells.Item(92, 3).Value = "displayLabel"
$trans.Cells.Item(92, 4).Value = "Left"
$trans.Cells.Item(92, 5).Value = "LTR"
$trans.Cells.Item(92, 6).Value = "IN<value>: <value>"

# Row 93
$trans.Cells.Item(93, 2).Value = "SingleUseId108"
$trans.Cells.Item(93, 3).Value = "displayLabel"
$trans.Cells.Item(93, 4).Value = "Center"
$trans.Cells.Item(93, 5).Value = "LTR"
$trans.Cells.Item(93, 6).Value = "<value>"

# Row 94
$trans.Cells.Item(94, 2).Value = "SingleUseId109"
$trans.Cells.Item(94, 3).Value = "displayMeas"
$trans.Cells.Item(94, 4).Value = "Left"
$trans.Cells.Item(94, 5).Value = "LTR"
$trans.Cells.Item(94, 6).Value = "Mean"

# Row 95
$trans.Cells.Item(95, 2).Value = "SingleUseId110"
$trans.Cells.Item(95, 3).Value = "displayMeas"
$trans.Cells.Item(95, 4).Value = "Left"
$trans.Cells.Item(95, 5).Value = "LTR"
$trans.Cells.Item(95, 6).Value = "StdDev"

# Row 96
$trans.Cells.Item(96, 2).Value = "SingleUseId111"
$trans.Cells.Item(96, 3).Value = "displayLabel"
$trans.Cells.Item(96, 4).Value = "Left"
$trans.Cells.Item(96, 5).Value = "LTR"
$trans.Cells.Item(96, 6).Value = "IN<value>: <value>"

# Row 97
$trans.Cells.Item(97, 2).Value = "SingleUseId112"
$trans.Cells.Item(97, 3).Value = "displayLabel"
$trans.Cells.Item(97, 4).Value = "Left"
$trans.Cells.Item(97, 5).Value = "LTR"
$trans.Cells.Item(97, 6).Value = "IN<value>: <value>"

# Row 98
$trans.Cells.Item(98, 2).Value = "SingleUseId113"
$trans.Cells.Item(98, 3).Value = "displayLabel"
$trans.Cells.Item(98, 4).Value = "Left"
$trans.Cells.Item(98, 5).Value = "LTR"
$trans.Cells.Item(98, 6).Value = "IN<value>: <value>"

# Row 99
$trans.Cells.Item(99, 2).Value = "SingleUseId114"
$trans.Cells.Item(99, 3).Value = "displayLabel"
$trans.Cells.Item(99, 4).Value = "Left"
$trans.Cells.Item(99, 5).Value = "LTR"
$trans.Cells.Item(99, 6).Value = "IN<value>: <value>"

# Row 100
$trans.Cells.Item(100, 2).Value = "SingleUseId115"
$trans.Cells.Item(100, 3).Value = "displayLabel"
$trans.Cells.Item(100, 4).Value = "Left"
$trans.Cells.Item(100, 5).Value = "LTR"
$trans.Cells.Item(100, 6).Value = "IN<value>: <value>"

# Row 101
$trans.Cells.Item(101, 2).Value = "SingleUseId116"
$trans.Cells.Item(101, 3).Value = "displayLabel"
$trans.Cells.Item(101, 4).Value = "Left"
$trans.Cells.Item(101, 5).Value = "LTR"
$trans.Cells.Item(101, 6).Value = "IN<value>: <value>"

# Row 102
$trans.Cells.Item(102, 2).Value = "SingleUseId117"
$trans.Cells.Item(102, 3).Value = "displayLabel"
$trans.Cells.Item(102, 4).Value = "Left"
$trans.Cells.Item(102, 5).Value = "LTR"
$trans.Cells.Item(102, 6).Value = "IN<value>: <value>"

# Row 103
$trans.Cells.Item(103, 2).Value = "SingleUseId118"
$trans.Cells.Item(103, 3).Value = "displayLabel"
$trans.Cells.Item(103, 4).Value = "Left"
$trans.Cells.Item(103, 5).Value = "LTR"
$trans.Cells.Item(103, 6).Value = "IN<value>: <value>"

# Row 104
$trans.Cells.Item(104, 2).Value = "SingleUseId119"
$trans.Cells.Item(104, 3).Value = "displayLabel"
$trans.Cells.Item(104, 4).Value = "Left"
$trans.Cells.Item(104, 5).Value = "LTR"
$trans.Cells.Item(104, 6).Value = "IN<value>: <value>"

# Row 105
$trans.Cells.Item(105, 2).Value = "SingleUseId120"
$trans.Cells.Item(105, 3).Value = "displayLabel"
$trans.Cells.Item(105, 4).Value = "Left"
$trans.Cells.Item(105, 5).Value = "LTR"
$trans.Cells.Item(105, 6).Value = "IN<value>: <value>"

# Row 106
$trans.Cells.Item(106, 2).Value = "SingleUseId121"
$trans.Cells.Item(106, 3).Value = "displayLabel"
$trans.Cells.Item(106, 4).Value = "Left"
$trans.Cells.Item(106, 5).Value = "LTR"
$trans.Cells.Item(106, 6).Value = "IN<value>: <value>"

# Row 107
$trans.Cells.Item(107, 2).Value = "SingleUseId122"
$trans.Cells.Item(107, 3).Value = "displayLabel"
$trans.Cells.Item(107, 4).Value = "Left"
$trans.Cells.Item(107, 5).Value = "LTR"
$trans.Cells.Item(107, 6).Value = "IN<value>: <value>"

# Row 108
$trans.Cells.Item(108, 2).Value = "SingleUseId123"
$trans.Cells.Item(108, 3).Value = "displayLabel"
$trans.Cells.Item(108, 4).Value = "Center"
$trans.Cells.Item(108, 5).Value = "LTR"
$trans.Cells.Item(108, 6).Value = "<value>"

# Row 109
$trans.Cells.Item(109, 2).Value = "SingleUseId124"
$trans.Cells.Item(109, 3).Value = "displayLabel"
$trans.Cells.Item(109, 4).Value = "Center"
$trans.Cells.Item(109, 5).Value = "LTR"
$trans.Cells.Item(109, 6).Value = "<value>"

# Row 110
$trans.Cells.Item(110, 2).Value = "SingleUseId125"
$trans.Cells.Item(110, 3).Value = "displayLabel"
$trans.Cells.Item(110, 4).Value = "Center"
$trans.Cells.Item(110, 5).Value = "LTR"
$trans.Cells.Item(110, 6).Value = "<value>"

# Row 111
$trans.Cells.Item(111, 2).Value = "SingleUseId126"
$trans.Cells.Item(111, 3).Value = "displayLabel"
$trans.Cells.Item(111, 4).Value = "Center"
$trans.Cells.Item(111, 5).Value = "LTR"
$trans.Cells.Item(111, 6).Value = "<value>"

# Row 112
$trans.Cells.Item(112, 2).Value = "SingleUseId127"
$trans.Cells.Item(112, 3).Value = "displayLabel"
$trans.Cells.Item(112, 4).Value = "Center"
$trans.Cells.Item(112, 5).Value = "LTR"
$trans.Cells.Item(112, 6).Value = "<value>"

# Row 113
$trans.Cells.Item(113, 2).Value = "SingleUseId128"
$trans.Cells.Item(113, 3).Value = "displayLabel"
$trans.Cells.Item(113, 4).Value = "Center"
$trans.Cells.Item(113, 5).Value = "LTR"
$trans.Cells.Item(113, 6).Value = "<value>"

# Row 114
$trans.Cells.Item(114, 2).Value = "SingleUseId129"
$trans.Cells.Item(114, 3).Value = "displayMeasValue"
$trans.Cells.Item(114, 4).Value = "Right"
$trans.Cells.Item(114, 5).Value = "LTR"
$trans.Cells.Item(114, 6).Value = "<value>"

# Row 115
$trans.Cells.Item(115, 2).Value = "SingleUseId130"
$trans.Cells.Item(115, 3).Value = "displayMeasValue"
$trans.Cells.Item(115, 4).Value = "Right"
$trans.Cells.Item(115, 5).Value = "LTR"
$trans.Cells.Item(115, 6).Value = "<value>"

# Row 116
$trans.Cells.Item(116, 2).Value = "SingleUseId131"
$trans.Cells.Item(116, 3).Value = "displayMeasValue"
$trans.Cells.Item(116, 4).Value = "Right"
$trans.Cells.Item(116, 5).Value = "LTR"
$trans.Cells.Item(116, 6).Value = "<value>"

# Row 117
$trans.Cells.Item(117, 2).Value = "SingleUseId132"
$trans.Cells.Item(117, 3).Value = "displayMeasValue"
$trans.Cells.Item(117, 4).Value = "Right"
$trans.Cells.Item(117, 5).Value = "LTR"
$trans.Cells.Item(117, 6).Value = "<value>"

# Row 118
$trans.Cells.Item(118, 2).Value = "SingleUseId133"
$trans.Cells.Item(118, 3).Value = "displayMeasValue"
$trans.Cells.Item(118, 4).Value = "Right"
$trans.Cells.Item(118, 5).Value = "LTR"
$trans.Cells.Item(118, 6).Value = "<value>"

# Row 119
$trans.Cells.Item(119, 2).Value = "SingleUseId134"
$trans.Cells.Item(119, 3).Value = "displayMeasValue"
$trans.Cells.Item(119, 4).Value = "Right"
$trans.Cells.Item(119, 5).Value = "LTR"
$trans.Cells.Item(119, 6).Value = "<value>"

# Row 120
$trans.Cells.Item(120, 2).Value = "SingleUseId135"
$trans.Cells.Item(120, 3).Value = "displayMeasValue"
$trans.Cells.Item(120, 4).Value = "Right"
$trans.Cells.Item(120, 5).Value = "LTR"
$trans.Cells.Item(120, 6).Value = "<value>"

# Row 121
$trans.Cells.Item(121, 2).Value = "SingleUseId136"
$trans.Cells.Item(121, 3).Value = "displayMeasValue"
$trans.Cells.Item(121, 4).Value = "Right"
$trans.Cells.Item(121, 5).Value = "LTR"
$trans.Cells.Item(121, 6).Value = "<value>"

# Row 122
$trans.Cells.Item(122, 2).Value = "SingleUseId137"
$trans.Cells.Item(122, 3).Value = "displayMeasValue"
$trans.Cells.Item(122, 4).Value = "Right"
$trans.Cells.Item(122, 5).Value = "LTR"
$trans.Cells.Item(122, 6).Value = "<value>"

# Row 123
$trans.Cells.Item(123, 2).Value = "SingleUseId138"
$trans.Cells.Item(123, 3).Value = "displayMeasValue"
$trans.Cells.Item(123, 4).Value = "Right"
$trans.Cells.Item(123, 5).Value = "LTR"
$trans.Cells.Item(123, 6).Value = "<value>"

# Row 124
$trans.Cells.Item(124, 2).Value = "SingleUseId139"
$trans.Cells.Item(124, 3).Value = "displayMeasValue"
$trans.Cells.Item(124, 4).Value = "Right"
$trans.Cells.Item(124, 5).Value = "LTR"
$trans.Cells.Item(124, 6).Value = "<value>"

# Row 125
$trans.Cells.Item(125, 2).Value = "SingleUseId140"
$trans.Cells.Item(125, 3).Value = "displayMeasValue"
$trans.Cells.Item(125, 4).Value = "Right"
$trans.Cells.Item(125, 5).Value = "LTR"
$trans.Cells.Item(125, 6).Value = "<value>"

# Row 126
$trans.Cells.Item(126, 2).Value = "SingleUseId141"
$trans.Cells.Item(126, 3).Value = "displayMeasValue"
$trans.Cells.Item(126, 4).Value = "Right"
$trans.Cells.Item(126, 5).Value = "LTR"
$trans.Cells.Item(126, 6).Value = "<value>"

# Row 127
$trans.Cells.Item(127, 2).Value = "SingleUseId142"
$trans.Cells.Item(127, 3).Value = "displayMeasValue"
$trans.Cells.Item(127, 4).Value = "Right"
$trans.Cells.Item(127, 5).Value = "LTR"
$trans.Cells.Item(127, 6).Value = "<value>"

Write-Output "Edit complete"